$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

# Measured sensor data (D:H) for rows 7-23 and 29-45, columns:
#   D = hm mozzo [mm], E = hr ref [mm], F = Laumas [N], G = Torsiometro [Nm], H = SG600 [V]
$ws.Range("D7").Value = 837
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -156.5579745837849
$ws.Range("G7").Value = -1.1772
$ws.Range("H7").Value = 0.001

$ws.Range("D8").Value = 837
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = -157.0411782090434
$ws.Range("G8").Value = -1.1772
$ws.Range("H8").Value = 0.001

$ws.Range("D9").Value = 837
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = -156.5579745837849
$ws.Range("G9").Value = -1.1772
$ws.Range("H9").Value = 0.001

$ws.Range("D10").Value = 837
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -157.0411782090434
$ws.Range("G10").Value = -1.1772
$ws.Range("H10").Value = 0.001

$ws.Range("D11").Value = 837
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -156.5579745837849
$ws.Range("G11").Value = -1.1772
$ws.Range("H11").Value = 0.001

$ws.Range("D12").Value = 837
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = -157.0411782090434
$ws.Range("G12").Value = -1.1772
$ws.Range("H12").Value = 0.001

$ws.Range("D13").Value = 837
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = -156.5579745837849
$ws.Range("G13").Value = -1.1772
$ws.Range("H13").Value = 0.001

$ws.Range("D14").Value = 837
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = -157.0411782090434
$ws.Range("G14").Value = -1.1772
$ws.Range("H14").Value = 0.001

$ws.Range("D15").Value = 837
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = -157.0411782090434
$ws.Range("G15").Value = -1.1772
$ws.Range("H15").Value = 0.001

$ws.Range("D16").Value = 837
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = -157.0411782090434
$ws.Range("G16").Value = -1.1772
$ws.Range("H16").Value = 0.001

$ws.Range("D17").Value = 837
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = -157.0411782090434
$ws.Range("G17").Value = -1.1772
$ws.Range("H17").Value = 0.001

$ws.Range("D18").Value = 837
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = -157.0411782090434
$ws.Range("G18").Value = -1.1772
$ws.Range("H18").Value = 0.001

$ws.Range("D19").Value = 837
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = -156.5579745837849
$ws.Range("G19").Value = -1.1772
$ws.Range("H19").Value = 0.001

$ws.Range("D20").Value = 837
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = -157.0411782090434
$ws.Range("G20").Value = -1.1772
$ws.Range("H20").Value = 0.001

$ws.Range("D21").Value = 837
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = -157.0411782090434
$ws.Range("G21").Value = -1.1772
$ws.Range("H21").Value = 0.001

$ws.Range("D22").Value = 837
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = -156.5579745837849
$ws.Range("G22").Value = -1.1772
$ws.Range("H22").Value = 0.001

$ws.Range("D23").Value = 837
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = -157.0411782090434
$ws.Range("G23").Value = -1.1772
$ws.Range("H23").Value = 0.001

$ws.Range("D29").Value = 837
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = -157.0411782090434
$ws.Range("G29").Value = -1.1772
$ws.Range("H29").Value = 0.001

$ws.Range("D30").Value = 837
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = -156.5579745837849
$ws.Range("G30").Value = -1.1772
$ws.Range("H30").Value = 0.001

$ws.Range("D31").Value = 837
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = -157.0411782090434
$ws.Range("G31").Value = -1.1772
$ws.Range("H31").Value = 0.001

$ws.Range("D32").Value = 837
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = -156.5579745837849
$ws.Range("G32").Value = -1.1772
$ws.Range("H32").Value = 0.001

$ws.Range("D33").Value = 837
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = -156.5579745837849
$ws.Range("G33").Value = -1.1772
$ws.Range("H33").Value = 0.001

$ws.Range("D34").Value = 837
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = -157.0411782090434
$ws.Range("G34").Value = -1.1772
$ws.Range("H34").Value = 0.001

$ws.Range("D35").Value = 837
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = -157.0411782090434
$ws.Range("G35").Value = -1.1772
$ws.Range("H35").Value = 0.001

$ws.Range("D36").Value = 837
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = -156.5579745837849
$ws.Range("G36").Value = -1.1772
$ws.Range("H36").Value = 0.001

$ws.Range("D37").Value = 837
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = -156.5579745837849
$ws.Range("G37").Value = -1.1772
$ws.Range("H37").Value = 0.001

$ws.Range("D38").Value = 837
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = -157.0411782090434
$ws.Range("G38").Value = -1.1772
$ws.Range("H38").Value = 0.001

$ws.Range("D39").Value = 837
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = -157.0411782090434
$ws.Range("G39").Value = -1.1772
$ws.Range("H39").Value = 0.001

$ws.Range("D40").Value = 837
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = -157.0411782090434
$ws.Range("G40").Value = -1.1772
$ws.Range("H40").Value = 0.001

$ws.Range("D41").Value = 837
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = -156.5579745837849
$ws.Range("G41").Value = -1.1772
$ws.Range("H41").Value = 0.001

$ws.Range("D42").Value = 837
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = -157.0411782090434
$ws.Range("G42").Value = -1.1772
$ws.Range("H42").Value = 0.001

$ws.Range("D43").Value = 837
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = -156.5579745837849
$ws.Range("G43").Value = -1.1772
$ws.Range("H43").Value = 0.001

$ws.Range("D44").Value = 837
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = -157.0411782090434
$ws.Range("G44").Value = -1.1772
$ws.Range("H44").Value = 0.001

$ws.Range("D45").Value = 837
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = -157.0411782090434
$ws.Range("G45").Value = -1.1772
$ws.Range("H45").Value = 0.001
